$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.458.77'
$ws.Range('E2').Value = '  +2.99%  '

$ws.Range('D3').Value = '3.768.86'
$ws.Range('E3').Value = '  +1.78%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.17'
$ws.Range('E5').Value = '  +1.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.26'
$ws.Range('E6').Value = '  +2.72%  '

$ws.Range('D7').Value = '3.767.23'
$ws.Range('E7').Value = '  +1.76%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('E9').Value = '  +2.11%  '

$ws.Range('E10').Value = '  +5.06%  '

$ws.Range('E11').Value = '  +3.26%  '

$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.55'
$ws.Range('E13').Value = '  +2.35%  '

$ws.Range('E14').Value = '  +4.56%  '

$ws.Range('D15').Value = '4.396.97'
$ws.Range('E15').Value = '  +1.81%  '

$ws.Range('D16').Value = '3.763.12'
$ws.Range('E16').Value = '  +1.66%  '

$ws.Range('D17').Value = '69.382.09'

$ws.Range('E18').Value = '  +2.13%  '

$ws.Range('E19').Value = '  +0.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.18'
$ws.Range('E20').Value = '  -2.12%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.92'
$ws.Range('E21').Value = '  +19.76%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '497.44'
$ws.Range('E22').Value = '  +1.14%  '

$ws.Range('E23').Value = '  +0.99%  '

$ws.Range('E24').Value = '  +12.26%  '

$ws.Range('E25').Value = '  -0.55%  '

$ws.Range('E26').Value = '  +2.14%  '

$ws.Range('E27').Value = '  +1.94%  '

$ws.Range('E28').Value = '  +2.29%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.00'
$ws.Range('E30').Value = '  +2.26%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.53'
$ws.Range('E31').Value = '  +7.25%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.00'
$ws.Range('E32').Value = '  +4.71%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.05'
$ws.Range('E33').Value = '  +1.71%  '

$ws.Range('D34').Value = '3.914.33'
$ws.Range('E34').Value = '  +2.02%  '

$ws.Range('E35').Value = '  +1.31%  '

$ws.Range('D36').Value = '3.700.99'
$ws.Range('E36').Value = '  +1.64%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('E38').Value = '  +2.03%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.93'
$ws.Range('E39').Value = '  +3.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.134'
$ws.Range('E40').Value = '  +1.83%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.327'
$ws.Range('E41').Value = '  +1.42%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.06'
$ws.Range('E42').Value = '  +9.76%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '440.51'
$ws.Range('E43').Value = '  +1.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.72'
$ws.Range('E44').Value = '  +0.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.00'
$ws.Range('E45').Value = '  +3.31%  '

$ws.Range('E46').Value = '  +1.32%  '

$ws.Range('E47').Value = '  +0.01%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.75'
$ws.Range('E48').Value = '  -0.09%  '

$ws.Range('D49').Value = '2.820.02'
$ws.Range('E49').Value = '  +2.46%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.15'
$ws.Range('E50').Value = '  -1.48%  '

$ws.Range('E51').Value = '  +2.53%  '
